$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.422.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.021.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +6.85%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.663'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.75%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '45.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.364'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0718'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.92%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.313.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.809'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.008.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.92'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.486.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.33%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.00%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0821'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.62%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '235.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.16%  '

$ws.Range("E23").Value = '  -6.08%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("E25").Value = '  -8.72%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.02%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.96%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -9.03%  '

$ws.Range("E30").Value = '  -4.07%  '

$ws.Range("B31").Value = 'Gas'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +56.51%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.39'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0592'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("E34").Value = '  -0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.03'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.59%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0823'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +13.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.842'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("E40").Value = '  -7.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0216'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.97%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.55%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.10'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +13.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.318.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.87%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0811'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.206.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.30%  '

$ws.Range("E51").Value = '  +14.12%  '
